$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 377.45456
$ws.Range("I9").Value = 88.71429000000001
$ws.Range("J9").Value = 882.75
$ws.Range("K9").Value = 88.71429000000001
$ws.Range("L9").Value = 882.75
$ws.Range("M9").Value = 80.28570999999999
$ws.Range("N9").Value = -1220.75
$ws.Range("H20").Value = 15266.667
$ws.Range("I20").Value = 15266.667
$ws.Range("K20").Value = 15266.667
$ws.Range("M20").Value = -15036.667
$ws.Range("H35").Value = 15266.667
$ws.Range("I35").Value = 15266.667
$ws.Range("K35").Value = 15266.667
$ws.Range("M35").Value = -14887.667
$ws.Range("H68").Value = 39800
$ws.Range("J68").Value = 39800
$ws.Range("L68").Value = 39800
$ws.Range("N68").Value = -41298
$ws.Range("H69").Value = 6000
$ws.Range("J69").Value = 6000
$ws.Range("L69").Value = 18000
$ws.Range("N69").Value = -19748
$ws.Range("H70").Value = 2313.5
$ws.Range("I70").Value = 3251.25
$ws.Range("J70").Value = 1375.75
$ws.Range("K70").Value = 9753.75
$ws.Range("L70").Value = 4127.25
$ws.Range("M70").Value = -9483.75
$ws.Range("N70").Value = -4667.25
$ws.Range("H71").Value = 39800
$ws.Range("J71").Value = 39800
$ws.Range("L71").Value = 119400
$ws.Range("N71").Value = -126888
$ws.Range("H72").Value = 6000
$ws.Range("J72").Value = 6000
$ws.Range("L72").Value = 54000
$ws.Range("N72").Value = -62736
$ws.Range("H73").Value = 2313.5
$ws.Range("I73").Value = 3251.25
$ws.Range("J73").Value = 1375.75
$ws.Range("K73").Value = 9753.75
$ws.Range("L73").Value = 4127.25
$ws.Range("M73").Value = -8817.75
$ws.Range("N73").Value = -5999.25
$ws.Range("H75").Value = 30001
$ws.Range("J75").Value = 30001
$ws.Range("L75").Value = 30001
$ws.Range("N75").Value = -31873
$ws.Range("H78").Value = 30001
$ws.Range("J78").Value = 30001
$ws.Range("L78").Value = 90003
$ws.Range("N78").Value = -99363
$ws.Range("H106").Value = 35715630
$ws.Range("I106").Value = 45455436
$ws.Range("K106").Value = 45455436
$ws.Range("M106").Value = -45454805
$ws.Range("H137").Value = 1852.6666
$ws.Range("I137").Value = 1299.8334
$ws.Range("J137").Value = 2405.5
$ws.Range("K137").Value = 3899.5002
$ws.Range("L137").Value = 7216.5
$ws.Range("M137").Value = -1349.5002
$ws.Range("N137").Value = -12316.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26673764
$ws.Range("I32").Value = 15157376
$ws.Range("J32").Value = 111127280
$ws.Range("K32").Value = 15157376
$ws.Range("L32").Value = 111127280
$ws.Range("M32").Value = -15157089
$ws.Range("N32").Value = -111127854
$ws.Range("H87").Value = 48333
$ws.Range("J87").Value = 48333
$ws.Range("L87").Value = 48333
$ws.Range("N87").Value = -50829
$ws.Range("H90").Value = 48333
$ws.Range("J90").Value = 48333
$ws.Range("L90").Value = 144999
$ws.Range("N90").Value = -157479

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 696.8125
$ws.Range("I94").Value = 666.5833
$ws.Range("J94").Value = 787.5
$ws.Range("K94").Value = 666.5833
$ws.Range("L94").Value = 787.5
$ws.Range("M94").Value = -215.5833
$ws.Range("N94").Value = -1689.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4232.9
$ws.Range("I62").Value = 2974.75
$ws.Range("J62").Value = 5071.6665
$ws.Range("K62").Value = 2974.75
$ws.Range("L62").Value = 5071.6665
$ws.Range("M62").Value = -2350.75
$ws.Range("N62").Value = -6319.6665
$ws.Range("H65").Value = 4232.9
$ws.Range("I65").Value = 2974.75
$ws.Range("J65").Value = 5071.6665
$ws.Range("K65").Value = 14873.75
$ws.Range("L65").Value = 25358.3325
$ws.Range("M65").Value = -11753.75
$ws.Range("N65").Value = -31598.3325
$ws.Range("H92").Value = 49900
$ws.Range("J92").Value = 49900
$ws.Range("L92").Value = 49900
$ws.Range("N92").Value = -54892
$ws.Range("H105").Value = 6683.3335
$ws.Range("I105").Value = 7594.2856
$ws.Range("K105").Value = 7594.2856
$ws.Range("M105").Value = -5847.2856
$ws.Range("H107").Value = 338.5366
$ws.Range("I107").Value = 237.5
$ws.Range("J107").Value = 614.0909
$ws.Range("K107").Value = 237.5
$ws.Range("L107").Value = 614.0909
$ws.Range("M107").Value = 1682.5
$ws.Range("N107").Value = -4454.0909

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 15000
$ws.Range("J75").Value = 15000
$ws.Range("L75").Value = 15000
$ws.Range("N75").Value = -16748
$ws.Range("H78").Value = 15000
$ws.Range("J78").Value = 15000
$ws.Range("L78").Value = 45000
$ws.Range("N78").Value = -53736
$ws.Range("H107").Value = 4730.0435
$ws.Range("I107").Value = 421.25
$ws.Range("J107").Value = 7028.067
$ws.Range("K107").Value = 421.25
$ws.Range("L107").Value = 7028.067
$ws.Range("M107").Value = 1498.75
$ws.Range("N107").Value = -10868.067
$ws.Range("H110").Value = 40949.5
$ws.Range("J110").Value = 40949.5
$ws.Range("L110").Value = 40949.5
$ws.Range("N110").Value = -49129.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 972.2857
$ws.Range("I46").Value = 302.66666
$ws.Range("J46").Value = 4990
$ws.Range("K46").Value = 302.66666
$ws.Range("L46").Value = 4990
$ws.Range("M46").Value = -114.66666
$ws.Range("N46").Value = -5366
$ws.Range("H136").Value = 3176.6072
$ws.Range("I136").Value = 2805.5
$ws.Range("J136").Value = 4289.9287
$ws.Range("K136").Value = 8416.5
$ws.Range("L136").Value = 12869.7861
$ws.Range("M136").Value = -5866.5
$ws.Range("N136").Value = -17969.7861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 29999.5
$ws.Range("J75").Value = 29999.5
$ws.Range("L75").Value = 29999.5
$ws.Range("N75").Value = -31871.5
$ws.Range("H78").Value = 29999.5
$ws.Range("J78").Value = 29999.5
$ws.Range("L78").Value = 89998.5
$ws.Range("N78").Value = -99358.5
$ws.Range("H100").Value = 721.0909
$ws.Range("I100").Value = 772.44446
$ws.Range("J100").Value = 490
$ws.Range("K100").Value = 1544.88892
$ws.Range("L100").Value = 980
$ws.Range("M100").Value = -1003.88892
$ws.Range("N100").Value = -2062
$ws.Range("H113").Value = 687.3333
$ws.Range("I113").Value = 655.4
$ws.Range("J113").Value = 751.2
$ws.Range("K113").Value = 1966.2
$ws.Range("L113").Value = 2253.6
$ws.Range("M113").Value = 203.8000000000002
$ws.Range("N113").Value = -6593.6
$ws.Range("H126").Value = 596.2889
$ws.Range("J126").Value = 1136.1428
$ws.Range("L126").Value = 3408.4284
$ws.Range("N126").Value = -8348.428400000001
$ws.Range("H136").Value = 2937.9824
$ws.Range("I136").Value = 3412.5
$ws.Range("J136").Value = 1988.9474
$ws.Range("K136").Value = 10237.5
$ws.Range("L136").Value = 5966.8422
$ws.Range("M136").Value = -7687.5
$ws.Range("N136").Value = -11066.8422
